# Re-organized item list in two functions for alignment with .csv and .xlsx
# config files: the placeholder database name in column D ("gailzsqlpool")
# is replaced throughout with "yourSynapseSQLPoolDbName", and the active
# selection on the sheet moves from I27:J27 to D3:D27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D (rows 2-27) holds the Synapse SQL pool database-name placeholder.
# Replace every occurrence so the now-unused "gailzsqlpool" shared string is
# dropped and a fresh "yourSynapseSQLPoolDbName" string is used instead.
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 4).Value2 = "yourSynapseSQLPoolDbName"
}

# Rows 25 and 26 previously carried an explicit (no-op) fill style on column D;
# drop it back to the workbook's Normal style so only those two cells lose
# their style index while the rest of the row keeps its formatting.
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"

# Move the active selection to D3:D27 (previously I27:J27).
[void]$ws.Range("D3:D27").Select()
